# Updates cryptocurrency price/volume data in the worksheet to reflect
# the latest GitHub Actions scrape (commit: "Updated cryptos list on
# Sat Jul 13 09:22:31 UTC 2024 with GitHub Actions").
#
# Only the "Price" (column D) and "Volume(1h)" (column E) cells that
# changed between scrapes are updated; everything else (coin name,
# link, rank) is left untouched. Price values that look like plain
# numbers are written with a leading apostrophe so Excel keeps them
# as text (matching how they were already stored in the workbook)
# instead of silently converting them to numeric cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.034.89"
$ws.Range("E2").Value = "  +1.49%  "
$ws.Range("D3").Value = "3.140.26"
$ws.Range("E3").Value = "  +2.14%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'536.24"
$ws.Range("E5").Value = "  +2.67%  "
$ws.Range("D6").Value = "'139.20"
$ws.Range("E6").Value = "  +2.68%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.512"
$ws.Range("E8").Value = "  +9.23%  "
$ws.Range("D9").Value = "'7.30"
$ws.Range("E9").Value = "  +0.33%  "
$ws.Range("E10").Value = "  +2.84%  "
$ws.Range("D11").Value = "'0.422"
$ws.Range("E11").Value = "  +5.42%  "
$ws.Range("E12").Value = "  +3.03%  "
$ws.Range("D13").Value = "3.675.11"
$ws.Range("E13").Value = "  +1.96%  "
$ws.Range("D14").Value = "'25.96"
$ws.Range("E14").Value = "  +3.22%  "
$ws.Range("D15").Value = "'0.0000169"
$ws.Range("E15").Value = "  +5.19%  "
$ws.Range("D16").Value = "58.110.39"
$ws.Range("E16").Value = "  +1.48%  "
$ws.Range("E17").Value = "  +6.18%  "
$ws.Range("D18").Value = "3.138.30"
$ws.Range("E18").Value = "  +2.21%  "
$ws.Range("D19").Value = "'12.98"
$ws.Range("E19").Value = "  +4.67%  "
$ws.Range("E20").Value = "  +4.40%  "
$ws.Range("D21").Value = "'375.20"
$ws.Range("E21").Value = "  +7.52%  "
$ws.Range("D22").Value = "'0.997"
$ws.Range("E22").Value = "  -0.30%  "
$ws.Range("E23").Value = "  -0.82%  "
$ws.Range("D24").Value = "'70.36"
$ws.Range("E24").Value = "  +2.08%  "
$ws.Range("D25").Value = "'0.516"
$ws.Range("E25").Value = "  +3.72%  "
$ws.Range("E26").Value = "  +1.31%  "
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("D28").Value = "'8.03"
$ws.Range("E28").Value = "  +12.04%  "
$ws.Range("D29").Value = "0.0₃0882"
$ws.Range("E29").Value = "  +2.88%  "
$ws.Range("D30").Value = "'1.90"
$ws.Range("E30").Value = "  +2.18%  "
$ws.Range("E31").Value = "  +6.86%  "
$ws.Range("D32").Value = "'21.75"
$ws.Range("E32").Value = "  +4.13%  "
$ws.Range("E33").Value = "  +7.00%  "
$ws.Range("E34").Value = "  +3.95%  "
$ws.Range("D35").Value = "'161.50"
$ws.Range("E35").Value = "  +1.42%  "
$ws.Range("D36").Value = "'6.24"
$ws.Range("E36").Value = "  +4.12%  "
$ws.Range("E37").Value = "  +9.75%  "
$ws.Range("D38").Value = "'25.42"
$ws.Range("E38").Value = "  +0.35%  "
$ws.Range("E39").Value = "  +7.80%  "
$ws.Range("D40").Value = "2.646.00"
$ws.Range("E40").Value = "  +9.75%  "
$ws.Range("D41").Value = "'0.0678"
$ws.Range("E41").Value = "  +3.57%  "
$ws.Range("D42").Value = "'4.25"
$ws.Range("E42").Value = "  +5.33%  "
$ws.Range("D43").Value = "'38.39"
$ws.Range("E43").Value = "  +5.28%  "
$ws.Range("D44").Value = "'0.700"
$ws.Range("E44").Value = "  +1.13%  "
$ws.Range("D45").Value = "'0.0275"
$ws.Range("E45").Value = "  +5.27%  "
$ws.Range("D46").Value = "'0.999"
$ws.Range("D47").Value = "'0.103"
$ws.Range("E47").Value = "  +12.55%  "
$ws.Range("D48").Value = "'6.22"
$ws.Range("E48").Value = "  +4.39%  "
$ws.Range("D49").Value = "'0.974"
$ws.Range("E49").Value = "  +3.94%  "
$ws.Range("D50").Value = "'20.21"
$ws.Range("E50").Value = "  +3.51%  "
$ws.Range("D51").Value = "'0.751"
$ws.Range("E51").Value = "  +0.24%  "
